# Add season-record columns (Wins, Losses, Ties) to the DET_2014 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold, border, centered) from an existing
# header cell onto the three new header cells, then set their text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (90 wins, 72 losses, 0 ties) for every player row.
for ($r = 2; $r -le 50; $r++) {
    $ws.Range("AD$r").Value = 90
    $ws.Range("AE$r").Value = 72
    $ws.Range("AF$r").Value = 0
}
